# Generate Report for Handoff
# Updates the localization-status report to reflect that the handoff
# package has been generated: status moves from "In Translation" to
# "Ready for handoff" and the related timestamps are refreshed. The
# "Status"/"zh-cn"/"de-de" columns are widened slightly to fit the new,
# longer status text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
# Closest width (in Excel "characters") achievable through the
# ColumnWidth API that rounds to the target stored column width.
$newStatusColWidth = 16.33

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-27 06:36:28"

$overview.Range("E1").ColumnWidth = $newStatusColWidth
$overview.Range("F1").ColumnWidth = $newStatusColWidth

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-27 06:36:24"

$zhcn.Range("C1").ColumnWidth = $newStatusColWidth

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-27 06:36:28"

$dede.Range("C1").ColumnWidth = $newStatusColWidth
